$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.066.42"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.870.79"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'598.90"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'167.76"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").Value = "3.870.87"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "'36.98"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "4.514.00"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "3.865.74"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "68.046.92"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "'18.21"
$ws.Range("E18").Value = "  +5.97%  "
$ws.Range("D19").Value = "'7.39"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'10.86"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'465.96"
$ws.Range("E22").Value = "  -4.26%  "
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("E24").Value = "  -4.75%  "
$ws.Range("D25").Value = "'83.32"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "'12.11"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "4.017.53"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "'7.74"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "'2.32"
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("D34").Value = "'31.21"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").Value = "'9.39"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("D36").Value = "3.840.75"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").Value = "'3.39"
$ws.Range("E38").Value = "  +7.41%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'0.312"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").Value = "'429.27"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "'1.98"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'47.23"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "'8.52"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'143.97"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000273"
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'40.44"
$ws.Range("E51").Value = "  +3.03%  "
